$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:34 PM"

# --- Sheet: Top Gainers ---
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsGainers.Range("B63").Value = "NPST"
$wsGainers.Range("C63").Value = 3.8509
$wsGainers.Range("D63").Value = -2.0059
$wsGainers.Range("E63").Value = -3.5057
$wsGainers.Range("B64").Value = "ORIENTTECH"
$wsGainers.Range("C64").Value = 3.827
$wsGainers.Range("D64").Value = 0.5247000000000001
$wsGainers.Range("E64").Value = 32.6784
$wsGainers.Range("B65").Value = "ICRA"
$wsGainers.Range("C65").Value = 3.7985
$wsGainers.Range("D65").Value = 4.4793
$wsGainers.Range("E65").Value = 2.8828
$wsGainers.Range("B66").Value = "SALASAR"
$wsGainers.Range("C66").Value = 3.7935
$wsGainers.Range("D66").Value = 4.7872
$wsGainers.Range("E66").Value = 11.0485

# --- Sheet: 1 Month Performance ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("B10").Value = "PANACHE"
$wsPerf.Range("C10").Value = 62.6487
$wsPerf.Range("B11").Value = "MAHASTEEL"
$wsPerf.Range("C11").Value = 55.9703
$wsPerf.Range("B12").Value = "INOXGREEN"
$wsPerf.Range("C12").Value = 51.0181
$wsPerf.Range("B13").Value = "STALLION"
$wsPerf.Range("C13").Value = 46.4325
$wsPerf.Range("B14").Value = "ORIENTTECH"
$wsPerf.Range("C14").Value = 45.3321
$wsPerf.Range("B15").Value = "TVSSRICHAK"
$wsPerf.Range("C15").Value = 40.7778
$wsPerf.Range("B16").Value = "MTARTECH"
$wsPerf.Range("C16").Value = 40.7213
$wsPerf.Range("B17").Value = "SEJALLTD"
$wsPerf.Range("C17").Value = 37.4301
$wsPerf.Range("B18").Value = "V2RETAIL"
$wsPerf.Range("C18").Value = 37.2004
$wsPerf.Range("B19").Value = "RAMAPHO"
$wsPerf.Range("C19").Value = 36.9731
$wsPerf.Range("B20").Value = "SANDUMA"
$wsPerf.Range("C20").Value = 36.9057
$wsPerf.Range("B21").Value = "TARACHAND"
$wsPerf.Range("C21").Value = 36.4813
$wsPerf.Range("B22").Value = "NETWEB"
$wsPerf.Range("C22").Value = 36.1199
$wsPerf.Range("B23").Value = "SAMMAANCAP"
$wsPerf.Range("C23").Value = 35.5128
$wsPerf.Range("B24").Value = "ONMOBILE"
$wsPerf.Range("C24").Value = 35.4702
$wsPerf.Range("B25").Value = "SHAREINDIA"
$wsPerf.Range("C25").Value = 35.3207
$wsPerf.Range("B26").Value = "SOUTHBANK"
$wsPerf.Range("C26").Value = 35.2819
$wsPerf.Range("B27").Value = "TVSELECT"
$wsPerf.Range("C27").Value = 35.1983
$wsPerf.Range("B28").Value = "RAMCOSYS"
$wsPerf.Range("C28").Value = 34.6928
$wsPerf.Range("B29").Value = "MAANALU"
$wsPerf.Range("C29").Value = 34.4803
$wsPerf.Range("B30").Value = "MEGASOFT"
$wsPerf.Range("C30").Value = 33.4399
$wsPerf.Range("B31").Value = "BHARATSE"
$wsPerf.Range("C31").Value = 31.8611
$wsPerf.Range("B32").Value = "EMKAY"
$wsPerf.Range("C32").Value = 30.3743
$wsPerf.Range("B33").Value = "ATHERENERG"
$wsPerf.Range("C33").Value = 29.116
$wsPerf.Range("B34").Value = "TATVA"
$wsPerf.Range("C34").Value = 28.7451
$wsPerf.Range("B35").Value = "TERASOFT"
$wsPerf.Range("C35").Value = 28.3093
$wsPerf.Range("B36").Value = "CARTRADE"
$wsPerf.Range("C36").Value = 27.5713
$wsPerf.Range("B37").Value = "ARFIN"
$wsPerf.Range("C37").Value = 27.3801
$wsPerf.Range("B38").Value = "MINDTECK"
$wsPerf.Range("C38").Value = 26.9415
$wsPerf.Range("B39").Value = "BHARATWIRE"
$wsPerf.Range("C39").Value = 26.5276
$wsPerf.Range("B40").Value = "HATSUN"
$wsPerf.Range("C40").Value = 26.492
$wsPerf.Range("B41").Value = "INDORAMA"
$wsPerf.Range("C41").Value = 26.4516
$wsPerf.Range("B42").Value = "IFBIND"
$wsPerf.Range("C42").Value = 26.161
$wsPerf.Range("B43").Value = "ADANIPOWER"
$wsPerf.Range("C43").Value = 25.8247
$wsPerf.Range("B44").Value = "AVALON"
$wsPerf.Range("C44").Value = 25.7352
$wsPerf.Range("B45").Value = "MRPL"
$wsPerf.Range("C45").Value = 25.6265
$wsPerf.Range("B46").Value = "HINDCOPPER"
$wsPerf.Range("C46").Value = 25.3164
$wsPerf.Range("B47").Value = "PRECWIRE"
$wsPerf.Range("C47").Value = 24.679
$wsPerf.Range("B48").Value = "SCI"
$wsPerf.Range("C48").Value = 24.132
$wsPerf.Range("B49").Value = "KICL"
$wsPerf.Range("C49").Value = 24.1119
$wsPerf.Range("B50").Value = "SKYGOLD"
$wsPerf.Range("C50").Value = 24.1079
$wsPerf.Range("B51").Value = "DCBBANK"
$wsPerf.Range("C51").Value = 23.8922
$wsPerf.Range("B52").Value = "AUBANK"
$wsPerf.Range("C52").Value = 23.6964
$wsPerf.Range("B53").Value = "ETHOSLTD"
$wsPerf.Range("C53").Value = 23.1527
$wsPerf.Range("B54").Value = "PVP"
$wsPerf.Range("C54").Value = 22.7524
$wsPerf.Range("B55").Value = "INDIANB"
$wsPerf.Range("C55").Value = 22.6689
$wsPerf.Range("B56").Value = "PRIVISCL"
$wsPerf.Range("C56").Value = 22.3984
$wsPerf.Range("B57").Value = "CPEDU"
$wsPerf.Range("C57").Value = 22.3786
$wsPerf.Range("B58").Value = "LORDSCHLO"
$wsPerf.Range("C58").Value = 22.1791
$wsPerf.Range("B59").Value = "GUJTHEM"
$wsPerf.Range("C59").Value = 22.0704
$wsPerf.Range("B60").Value = "SURYODAY"
$wsPerf.Range("C60").Value = 21.8039
$wsPerf.Range("B61").Value = "TDPOWERSYS"
$wsPerf.Range("C61").Value = 21.7743
$wsPerf.Range("B62").Value = "ORBTEXP"
$wsPerf.Range("C62").Value = 21.6115
$wsPerf.Range("B63").Value = "CEATLTD"
$wsPerf.Range("C63").Value = 20.0239
$wsPerf.Range("B64").Value = "ATL"
$wsPerf.Range("C64").Value = 19.9362
$wsPerf.Range("B65").Value = "GRMOVER"
$wsPerf.Range("C65").Value = 19.7859
$wsPerf.Range("B66").Value = "FEDERALBNK"
$wsPerf.Range("C66").Value = 19.6872
$wsPerf.Range("B67").Value = "SUBROS"
$wsPerf.Range("C67").Value = 19.6508
$wsPerf.Range("B68").Value = "USHAMART"
$wsPerf.Range("C68").Value = 19.6172
$wsPerf.Range("B69").Value = "BANKINDIA"
$wsPerf.Range("C69").Value = 19.3067
$wsPerf.Range("B70").Value = "RBLBANK"
$wsPerf.Range("C70").Value = 19.2556
$wsPerf.Range("B71").Value = "MOLDTECH"
$wsPerf.Range("C71").Value = 19.1891
$wsPerf.Range("B72").Value = "THOMASCOTT"
$wsPerf.Range("C72").Value = 19.1649
$wsPerf.Range("B74").Value = "KARURVYSYA"
$wsPerf.Range("C74").Value = 18.8614
$wsPerf.Range("B75").Value = "LUMAXIND"
$wsPerf.Range("C75").Value = 18.8057
$wsPerf.Range("B76").Value = "REPRO"
$wsPerf.Range("C76").Value = 18.689
